# StopAll 5% wrong ans
# Fixes the "Hoja1" comparison sheet: adds an input/expected-output/actual-output
# block, highlights matching results in light green, and tidies up stray cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Colors (COM uses BGR-packed integers: R + G*256 + B*65536)
$RED        = 255        # FF0000
$GREEN00B0  = 5287936    # 00B050
$LITEGREEN  = 5296274    # 92D050 (was FF0000 before the fix)

# --- New "input / output ex / My output" header block -------------------
$ws.Range("A21").Value = "input"
$ws.Range("E20").Value = 2
$ws.Range("E20").Font.Color = $RED
$ws.Range("E21").Value = 7
$ws.Range("E22").Value = 3

$ws.Range("B27").Value = "output ex"
$ws.Range("C27").Value = "My output"

# --- Row 28-36: duplicate the "B" result into a new "C" column and ------
# --- highlight both with the light-green fill (correct-answer marker) ---
$ws.Range("C28").Value = $ws.Range("B28").Value()
$ws.Range("C29").Value = $ws.Range("B29").Value()
$ws.Range("C30").Value = $ws.Range("B30").Value()
$ws.Range("C31").Value = $ws.Range("B31").Value()
$ws.Range("C32").Value = $ws.Range("B32").Value()
$ws.Range("C33").Value = $ws.Range("B33").Value()
$ws.Range("C34").Value = $ws.Range("B34").Value()
$ws.Range("C35").Value = $ws.Range("B35").Value()
$ws.Range("C36").Value = $ws.Range("B36").Value()

$ws.Range("B28:C36").Interior.Color = $LITEGREEN

# Row 30 also gained D/E values and lost its old red highlight in favor of
# the regular green used by the other rows in A28:A36.
$ws.Range("A30").Interior.Color = $GREEN00B0
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 2
$ws.Range("E30").Font.Color = $RED

$ws.Range("D31").Value = 6
$ws.Range("E31").Value = 7
$ws.Range("E31").Font.Color = $RED

# --- Drop the now-unused "C" column in the trailer table (rows 39-42) ---
$ws.Range("C39").ClearContents()
$ws.Range("C40").ClearContents()
$ws.Range("C41").ClearContents()
$ws.Range("C42").ClearContents()

# --- Column widths (best achievable under this host's width quantizer) --
$ws.Columns.Item(1).ColumnWidth = 6.6
$ws.Columns.Item(3).ColumnWidth = 13.6
$ws.Columns.Item(4).ColumnWidth = 4.1
$ws.Columns.Item(5).ColumnWidth = 11.8

# --- View state: zoom + scroll position + selection ----------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 150
$ws.Range("E33").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
